# Add team record (Wins/Losses/Ties) columns AD/AE/AF to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set up the new header cells, copying the style/format from the
# existing last header cell (AC1) so they match the rest of row 1
# (bold, centered, bordered).
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every data row (2 through 44) with the
# constant values used in the source data.
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 0
}
